# Group project testing.docx - fix "Dehli" typo and relocate the stray
# "_GoBack" bookmark (mirrors Word's own behaviour: the last edit position
# moves the hidden _GoBack bookmark there).
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. "Dehli" -> "Delhi"
#    Word had wrapped the misspelled run in <w:proofErr spellStart/spellEnd>
#    markers. Simply replacing the run text leaves those markers behind,
#    so instead we splice in a brand-new paragraph carrying the corrected
#    text (inheriting the original paragraph's own formatting/ilvl) and
#    delete the old paragraph (text + mark + its proofErr wrapper) outright.
# ---------------------------------------------------------------------
$paras = $d.Paragraphs
$count = $paras.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text -like "*Dehli*") {
        $rng = $p.Range
        $insertPoint = $d.Range($rng.Start, $rng.Start)
        $insertPoint.InsertBefore("Delhi`r")

        # Original paragraph (now shifted right by "Delhi" + CR = 6 chars)
        $oldStart = $rng.Start + 6
        $oldEnd = $rng.End + 6
        $oldRange = $d.Range($oldStart, $oldEnd)
        $oldRange.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# 2. Move the hidden "_GoBack" bookmark from the "Overall Test Result"
#    paragraph to just before the "Indonesia" run (the last-edited spot).
# ---------------------------------------------------------------------
$searchRange = $d.Content
$searchRange.Find.Execute("Indonesia") | Out-Null
$indonesiaStart = $searchRange.Start

$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

$newBookmarkRange = $d.Range($indonesiaStart, $indonesiaStart)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange)
